$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.258.16"
$ws.Range("E2").Value = "  -4.42%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.485.61"
$ws.Range("E3").Value = "  -3.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.69"
$ws.Range("E5").Value = "  -2.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.33"
$ws.Range("E6").Value = "  -6.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.575"
$ws.Range("E8").Value = "  -3.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.513.09"
$ws.Range("E9").Value = "  -2.84%  "

$ws.Range("E10").Value = "  -4.55%  "

$ws.Range("E11").Value = "  -2.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.59"
$ws.Range("E12").Value = "  +0.49%  "

$ws.Range("E13").Value = "  -3.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.952.44"
$ws.Range("E14").Value = "  -2.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.99"
$ws.Range("E15").Value = "  -6.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.197.71"
$ws.Range("E16").Value = "  -4.40%  "

$ws.Range("E17").Value = "  -3.87%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.506.01"
$ws.Range("E18").Value = "  -3.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.28"
$ws.Range("E19").Value = "  -3.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("E20").Value = "  -5.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.57"
$ws.Range("E21").Value = "  -4.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("E23").Value = "  -4.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.29"
$ws.Range("E24").Value = "  -3.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.440"
$ws.Range("E25").Value = "  -10.72%  "

$ws.Range("E26").Value = "  -3.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").Value = "  -0.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.608.25"
$ws.Range("E28").Value = "  -3.42%  "

$ws.Range("E29").Value = "  -4.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.87"
$ws.Range("E30").Value = "  -5.74%  "

$ws.Range("E31").Value = "  -7.11%  "

$ws.Range("E32").Value = "  -7.16%  "

$ws.Range("E33").Value = "  -5.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.71"
$ws.Range("E35").Value = "  -2.55%  "

$ws.Range("E36").Value = "  +2.06%  "

$ws.Range("E37").Value = "  -3.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.44"
$ws.Range("E38").Value = "  -9.15%  "

$ws.Range("E39").Value = "  -10.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.90"
$ws.Range("E40").Value = "  -2.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "307.05"
$ws.Range("E41").Value = "  -6.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.77"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.68"
$ws.Range("E43").Value = "  -6.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.822"
$ws.Range("E44").Value = "  -9.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.996"
$ws.Range("E45").Value = "  -0.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.595"
$ws.Range("E46").Value = "  -2.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.79"
$ws.Range("E47").Value = "  -1.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.19"
$ws.Range("E48").Value = "  +0.89%  "

$ws.Range("E49").Value = "  -3.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.61"
$ws.Range("E50").Value = "  -4.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0517"
$ws.Range("E51").Value = "  -5.75%  "
